# O4+ : add a new "Unc at peak" column.
# The uncertainty (col C) relative to the cross section (col B) at its peak
# value, i.e. the cross-section uncertainty expressed as a fraction at the
# row where the cross section (B:B) is maximal.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F1").Value = "Unc at peak"
$ws.Range("F2").Formula = "=XLOOKUP(MAX(B:B),B:B,C:C)/MAX(B:B)"
